# Apply crypto price/volume updates scraped on Wed Jul 31 19:33:44 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of "Price" cells in column D hold plain decimal strings (e.g. "0.640",
# "1.00"). Assigning those through .Value would make Excel auto-detect them as
# numbers and silently reformat/trim them (dropping trailing zeros, switching to
# scientific notation, etc). Mark those specific cells as Text first so the scraped
# price strings round-trip byte-for-byte, same as the rest of column D.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

$ws.Range('D2').Value = '65.745.15'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '3.276.98'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '585.63'
$ws.Range('E5').Value = '  +2.50%  '
$ws.Range('D6').Value = '177.79'
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('D7').Value = '0.640'
$ws.Range('E7').Value = '  +2.35%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '3.273.87'
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('D10').Value = '0.124'
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('E11').Value = '  +2.29%  '
$ws.Range('D12').Value = '0.398'
$ws.Range('E12').Value = '  +0.40%  '
$ws.Range('D13').Value = '3.860.44'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').Value = '0.129'
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('D15').Value = '65.851.69'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').Value = '26.36'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.0000162'
$ws.Range('E17').Value = '  +0.27%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.275.92'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').Value = '420.66'
$ws.Range('E19').Value = '  -3.40%  '
$ws.Range('D20').Value = '5.46'
$ws.Range('E20').Value = '  -1.27%  '
$ws.Range('D21').Value = '12.97'
$ws.Range('E21').Value = '  -0.98%  '
$ws.Range('D22').Value = '7.24'
$ws.Range('E22').Value = '  -1.75%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').Value = '70.96'
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('D25').Value = '5.65'
$ws.Range('E25').Value = '  -0.48%  '
$ws.Range('D26').Value = '0.206'
$ws.Range('E26').Value = '  +5.82%  '
$ws.Range('D27').Value = '0.506'
$ws.Range('E27').Value = '  +0.45%  '
$ws.Range('D28').Value = '0.0000113'
$ws.Range('E28').Value = '  +1.11%  '
$ws.Range('D29').Value = '9.37'
$ws.Range('E29').Value = '  +6.19%  '
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('D31').Value = '1.91'
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('D32').Value = '22.17'
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = '5.12'
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('D35').Value = '6.55'
$ws.Range('E35').Value = '  -0.24%  '
$ws.Range('D36').Value = '1.18'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Value = '157.85'
$ws.Range('E37').Value = '  -1.35%  '
$ws.Range('D38').Value = '1.43'
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('D39').Value = '2.848.48'
$ws.Range('E39').Value = '  +3.29%  '
$ws.Range('D40').Value = '1.78'
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('D41').Value = '26.09'
$ws.Range('E41').Value = '  -2.07%  '
$ws.Range('D42').Value = '4.31'
$ws.Range('E42').Value = '  +0.43%  '
$ws.Range('D43').Value = '0.746'
$ws.Range('E43').Value = '  -3.35%  '
$ws.Range('D44').Value = '39.53'
$ws.Range('E44').Value = '  -1.72%  '
$ws.Range('D45').Value = '5.87'
$ws.Range('E45').Value = '  -2.34%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '0.0635'
$ws.Range('E46').Value = '  -2.83%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = '2.28'
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('D48').Value = '310.65'
$ws.Range('E48').Value = '  -1.82%  '
$ws.Range('D49').Value = '22.76'
$ws.Range('E49').Value = '  -2.06%  '
$ws.Range('D50').Value = '0.0267'
$ws.Range('E50').Value = '  +0.58%  '
$ws.Range('E51').Value = '  -0.17%  '
